$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Marzo de 2020 a las 05:46"

# Refresh case counts for the countries around the 44-46 ranking band.
# The data refresh also causes India to overtake Catar and Filipinas,
# so the country names shift up while new totals are recorded.

# Row 44 (was Catar) -> now India
$ws.Range("A44").Value = "India"
$ws.Range("B44").Value = 504
$ws.Range("C44").Value = 5
$ws.Range("D44").Value = 37
$ws.Range("E44").Value = 457
$ws.Range("F44").Value = 0
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 10

# Row 45 (was Filipinas) -> now Catar
$ws.Range("A45").Value = "Catar"
$ws.Range("B45").Value = 501
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 33
$ws.Range("E45").Value = 468
$ws.Range("F45").Value = 6
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0

# Row 46 (was India) -> now Filipinas
$ws.Range("A46").Value = "Filipinas"
$ws.Range("B46").Value = 501
$ws.Range("C46").Value = 39
$ws.Range("D46").Value = 19
$ws.Range("E46").Value = 449
$ws.Range("F46").Value = 1
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 33
